$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-missing values for row 7 (C7, E7)
$ws.Range("C7").Value = 98.8
$ws.Range("E7").Value = 27.7

# Add new row 8 with the day's measurements
$ws.Range("A8").Value = 45946
$ws.Range("A8").NumberFormat = $ws.Range("A7").NumberFormat

$ws.Range("B8").Value = 98
$ws.Range("C8").Value = 98.09999999999999
$ws.Range("D8").Value = 28.6
$ws.Range("E8").Value = 27.4
